$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D. This shifts the existing D..AI
# columns one place to the right (D->E, E->F, ... AI->AJ), carrying along
# cell styles, merged cells, data validations and column widths.
$ws.Columns("D").Insert()

# The newly inserted column D gets a header in row 1 ("Criterios de
# aceptacion"); the data rows in column D stay blank (matching the source
# edit, which only adds the header label).
$ws.Range("D1").Value = "Criterios de aceptacion"

# Update the print area to include the extra column.
$ws.PageSetup.PrintArea = '$B$1:$N$27'

# Re-point the "CLICK HERE TO CREATE IN SMARTSHEET" hyperlink so its range
# keeps pace with the widened table (B29:M29 -> B29:N29).
$link = $ws.Hyperlinks.Item(1)
$address = $link.Address
$textToDisplay = $link.TextToDisplay
$link.Delete()
$ws.Hyperlinks.Add($ws.Range("B29:N29"), $address, "", "", $textToDisplay)

# Move the active selection to D5, matching the saved view state.
$ws.Range("D5").Select()
